$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update NUC_NASZ scenario figures in column E (rows 27-36) ---
# "Scenariusz NASZ" block (rows 27-31)
$ws.Range("E27").Value = 0
$ws.Range("E28").Value = 100
$ws.Range("E29").Value = 200
$ws.Range("E30").Value = 350
$ws.Range("E31").Value = 470

# second block (rows 32-36)
$ws.Range("E32").Value = 40
$ws.Range("E33").Value = 90
$ws.Range("E34").Value = 130
$ws.Range("E35").Value = 180
$ws.Range("E36").Value = 230

# --- Update the saved view state: scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H31").Select()
